# Upload excel files with prices
#
# The crawler re-ran and re-emitted the sheet: every row's "timestamp"
# (column O) is refreshed, and two pairs of product rows that used to be
# emitted together in one order now come out in the opposite order, so
# their whole row of data (columns A..N) swaps between row 17 <-> row 19
# and row 23 <-> row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-09-11 20:59:21"

# Write a single cell value, forcing Excel to store numeric-looking text
# (ids like "6577910", prices like "2.95") as text rather than silently
# re-typing it as a number - same trick a real user has to use when
# typing a leading-zero/price-like string into a "General" formatted
# cell. ClearFormats() afterwards drops the quote-prefix/text style we
# had to apply so the cell keeps its original (default) style.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

function Set-NumberValue($row, $col, $number) {
    $ws.Cells.Item($row, $col).Value = $number
}

function Clear-Value($row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

# --- row 17: was Kandoo (3630714) -> becomes Naturaline Kids (6577910) ---
Set-TextValue   17 1  "6577910"
$ws.Cells.Item(17, 2).Value  = "Naturaline Kids feuchtes Toilettenpapier"
$ws.Cells.Item(17, 3).Value  = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/naturaline-kids-feuchtes-toilettenpapier/p/6577910"
$ws.Cells.Item(17, 4).Value  = "50ST"
Set-NumberValue 17 5  2
Set-NumberValue 17 6  3
$ws.Cells.Item(17, 7).Value  = "Coop"
Set-TextValue   17 8  "2.95"
$ws.Cells.Item(17, 9).Value  = "0.06/1ST"
$ws.Cells.Item(17, 10).Value = "Preis pro 1 Stück"
Set-TextValue   17 11 "0.06"
$ws.Cells.Item(17, 12).Value = "1ST"
$ws.Cells.Item(17, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Cells.Item(17, 14).Value = "Naturaline Kids feuchtes Toilettenpapier 25% ab 2 Aktion 2.95 Schweizer Franken"

# --- row 19: was Naturaline Kids (6577910) -> becomes Kandoo (3630714) ---
Set-TextValue   19 1  "3630714"
$ws.Cells.Item(19, 2).Value  = "Kandoo Feuchttücher Melone 55 Stück"
$ws.Cells.Item(19, 3).Value  = "/de/kosmetik-gesundheit/baby-kind/pflege-accessoires/feuchttuecher/kandoo-feuchttuecher-melone-55-stueck/p/3630714"
$ws.Cells.Item(19, 4).Value  = "55ST"
Set-NumberValue 19 5  4
Set-NumberValue 19 6  5
$ws.Cells.Item(19, 7).Value  = "Kandoo"
Set-TextValue   19 8  "3.75"
$ws.Cells.Item(19, 9).Value  = "0.07/1ST"
$ws.Cells.Item(19, 10).Value = "Preis pro 1 Stück"
Set-TextValue   19 11 "0.07"
$ws.Cells.Item(19, 12).Value = "1ST"
$ws.Cells.Item(19, 13).Value = "['kosmetik-gesundheit', 'baby-kind', 'pflege-accessoires', 'feuchttuecher']"
$ws.Cells.Item(19, 14).Value = "Kandoo Feuchttücher Melone 55 Stück 3.75 Schweizer Franken"

# --- row 23: was Oecoplan (3874909) -> becomes Tela Viva (6996030) ---
Set-TextValue   23 1  "6996030"
$ws.Cells.Item(23, 2).Value  = "Tela Viva Haushaltspapier 3-lagig 4 Rollen"
$ws.Cells.Item(23, 3).Value  = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/tela-viva-haushaltspapier-3-lagig-4-rollen/p/6996030"
$ws.Cells.Item(23, 4).Value  = "200BLT"
Set-NumberValue 23 5  1
Set-NumberValue 23 6  4
$ws.Cells.Item(23, 7).Value  = "Tela"
Set-TextValue   23 8  "5.95"
Clear-Value     23 9
Clear-Value     23 10
Clear-Value     23 11
Clear-Value     23 12
$ws.Cells.Item(23, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Cells.Item(23, 14).Value = "Tela Viva Haushaltspapier 3-lagig 4 Rollen 5.95 Schweizer Franken"

# --- row 24: was Tela Viva (6996030) -> becomes Oecoplan (3874909) ---
Set-TextValue   24 1  "3874909"
$ws.Cells.Item(24, 2).Value  = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück"
$ws.Cells.Item(24, 3).Value  = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/oecoplan-papiertaschentuecher-special-edition-calendula-30x10-stueck/p/3874909"
$ws.Cells.Item(24, 4).Value  = "30ST"
Set-NumberValue 24 5  1
Set-NumberValue 24 6  5
$ws.Cells.Item(24, 7).Value  = "Coop"
Set-TextValue   24 8  "3.65"
$ws.Cells.Item(24, 9).Value  = "0.12/1ST"
$ws.Cells.Item(24, 10).Value = "Preis pro 1 Stück"
Set-TextValue   24 11 "0.12"
$ws.Cells.Item(24, 12).Value = "1ST"
$ws.Cells.Item(24, 13).Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Cells.Item(24, 14).Value = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück 20% Aktion 3.65 Schweizer Franken statt 4.60 Schweizer Franken"

# --- refresh the timestamp column for every data row ---
for ($row = 2; $row -le 30; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
